$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.288329
$ws.Range("H2").Value = 117.864987
$ws.Range("I2").Value = 0.632237668435316
$ws.Range("J2").Value = 0.632237668435316
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 649.656250165716
$ws.Range("R2").Value = 5846.906251491444
$ws.Range("S2").Value = 0.1340630353364149
$ws.Range("T2").Value = 0.1340630353364149
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 39.288329
$ws.Range("H3").Value = 117.864987
$ws.Range("I3").Value = 0.632237668435316
$ws.Range("J3").Value = 0.632237668435316
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 1595.916701819489
$ws.Range("R3").Value = 14363.2503163754
$ws.Range("S3").Value = 0.3293333007039109
$ws.Range("T3").Value = 0.3293333007039109
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 39.288329
$ws.Range("H4").Value = 117.864987
$ws.Range("I4").Value = 0.632237668435316
$ws.Range("J4").Value = 0.632237668435316
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 818.1884484523409
$ws.Range("R4").Value = 7363.696036071069
$ws.Range("S4").Value = 0.1688413323949903
$ws.Range("T4").Value = 0.1688413323949903
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.344283
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3112930657211948
$ws.Range("J5").Value = 0.3112930657211947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 319.869403351932
$ws.Range("R5").Value = 2878.824630167388
$ws.Range("S5").Value = 0.0660082360689509
$ws.Range("T5").Value = 0.06600823606895088
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.344283
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3112930657211948
$ws.Range("J6").Value = 0.3112930657211947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 785.7769752544788
$ws.Range("R6").Value = 7071.992777290308
$ws.Range("S6").Value = 0.1621529022051447
$ws.Range("T6").Value = 0.1621529022051447
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.344283
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3112930657211948
$ws.Range("J7").Value = 0.3112930657211947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 402.849123315807
$ws.Range("R7").Value = 3625.642109842263
$ws.Range("S7").Value = 0.08313192744709914
$ws.Range("T7").Value = 0.08313192744709913
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.509096666666667
$ws.Range("H8").Value = 10.52729
$ws.Range("I8").Value = 0.05646926584348937
$ws.Range("J8").Value = 0.05646926584348937
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 58.02503287772
$ws.Range("R8").Value = 522.22529589948
$ws.Range("S8").Value = 0.0119740432437895
$ws.Range("T8").Value = 0.0119740432437895
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.509096666666667
$ws.Range("H9").Value = 10.52729
$ws.Range("I9").Value = 0.05646926584348937
$ws.Range("J9").Value = 0.05646926584348937
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 142.5417196702978
$ws.Range("R9").Value = 1282.87547703268
$ws.Range("S9").Value = 0.02941490302940664
$ws.Range("T9").Value = 0.02941490302940664
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.509096666666667
$ws.Range("H10").Value = 10.52729
$ws.Range("I10").Value = 0.05646926584348937
$ws.Range("J10").Value = 0.05646926584348937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 73.07774166647
$ws.Range("R10").Value = 657.6996749982301
$ws.Range("S10").Value = 0.01508031957029324
$ws.Range("T10").Value = 0.01508031957029324
